# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Vega Modelo de Temuco - Brócoli" just
# after the existing row 636, pushing the old rows 637-697 down to 639-699.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 637:638; everything from the old row 637 onward
# shifts down by two rows (old 637 -> new 639, ..., old 697 -> new 699).
$ws.Rows("637:638").Insert()

# New row 637
$ws.Range("A637").Value = 10
$ws.Range("B637").Value = "Vega Modelo de Temuco"
$ws.Range("C637").Value = "La Araucanía"
$ws.Range("D637").Value = 45106
$ws.Range("E637").Value = 9
$ws.Range("F637").Value = 100112023
$ws.Range("G637").Value = "Brócoli"
$ws.Range("H637").Value = "Sin especificar"
$ws.Range("I637").Value = "Primera"
$ws.Range("J637").Value = 2800
$ws.Range("K637").Value = 1200
$ws.Range("L637").Value = 1200
$ws.Range("M637").Value = 1200
$ws.Range("N637").Value = "$/unidad"
$ws.Range("O637").Value = "Región Metropolitana"
$ws.Range("P637").Value = 1200
$ws.Range("Q637").Value = 1
$ws.Range("R637").Value = "Hortaliza"

# New row 638
$ws.Range("A638").Value = 10
$ws.Range("B638").Value = "Vega Modelo de Temuco"
$ws.Range("C638").Value = "La Araucanía"
$ws.Range("D638").Value = 45106
$ws.Range("E638").Value = 9
$ws.Range("F638").Value = 100112023
$ws.Range("G638").Value = "Brócoli"
$ws.Range("H638").Value = "Sin especificar"
$ws.Range("I638").Value = "Primera"
$ws.Range("J638").Value = 2150
$ws.Range("K638").Value = 1300
$ws.Range("L638").Value = 1300
$ws.Range("M638").Value = 1300
$ws.Range("N638").Value = "$/unidad"
$ws.Range("O638").Value = "Región de O'Higgins"
$ws.Range("P638").Value = 1300
$ws.Range("Q638").Value = 1
$ws.Range("R638").Value = "Hortaliza"
